$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "[-, Leonardo-Manut. Mecânica, -, Emerson-Robótica]"
$ws.Range("E2").Value = "[Nilton Maia-Elementos de máquinas, Leonardo-Manut. Mecânica, Emerson-Robótica, Ludoff-Eletrohidráulica]"

$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "[-, Leonardo-Manut. Mecânica, -, Euclides-Tecnologia da soldagem]"
$ws.Range("E3").Value = "[Nilton Maia-Elementos de máquinas, Victor Lima-CAM, Carlos Eduardo-Processos de Usinagem 1, Rogério-Processos de Usinagem 2]"
$ws.Range("F3").Value = "Cleidson-Máquinas Elétri"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "[Pedro Bispo-Automação Industrial, Humberto-Eletropneumática, Eudes-Microcontroladores, Euclides-Tecnologia da soldagem]"
$ws.Range("E4").Value = "[Nilton Maia-Elementos de máquinas, Victor Lima-CAM, Carlos Eduardo-Processos de Usinagem 1, Rogério-Processos de Usinagem 2]"
$ws.Range("F4").Value = "Cleidson-Máquinas Elétri"

$ws.Range("B6").Value = "-"
$ws.Range("D6").Value = "[Pedro Bispo-Automação Industrial, Humberto-Eletropneumática, Eudes-Microcontroladores, Euclides-Tecnologia da soldagem]"
$ws.Range("E6").Value = "[Nilton Maia-Elementos de máquinas, Victor Lima-CAM, Carlos Eduardo-Processos de Usinagem 1, Rogério-Processos de Usinagem 2]"
$ws.Range("F6").Value = "André Guimarães-Máquinas Térmicas e de Fl"

$ws.Range("B7").Value = "[-, Emerson-Robótica, -, -]"
$ws.Range("C7").Value = "[-, Pedro Bispo-Lab. Máquinas Elétricas, João Paulo-Lab. de eletroeletrônica]"
$ws.Range("D7").Value = "[Pedro Bispo-Automação Industrial, Humberto-Eletropneumática, Eudes-Microcontroladores, Ludoff-Eletrohidráulica]"
$ws.Range("E7").Value = "[Ludoff-Eletrohidráulica, Victor Lima-CAM, Carlos Eduardo-Processos de Usinagem 1, Rogério-Processos de Usinagem 2]"
$ws.Range("F7").Value = "André Guimarães-Máquinas Térmicas e de Fl"

$ws.Range("C8").Value = "[-, Pedro Bispo-Lab. Máquinas Elétricas, João Paulo-Lab. de eletroeletrônica]"
$ws.Range("D8").Value = "[Pedro Bispo-Automação Industrial, Humberto-Eletropneumática, Eudes-Microcontroladores, Euclides-Tecnologia da soldagem]"
$ws.Range("E8").Value = "[Ludoff-Eletrohidráulica, Emerson-Robótica, Leonardo-Manut. Mecânica, -]"
